$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 550
$ws.Range("I28").Value = 550
$ws.Range("K28").Value = 550
$ws.Range("M28").Value = -65

$ws.Range("H43").Value = 3438.1428
$ws.Range("I43").Value = 2953.4
$ws.Range("J43").Value = 4650
$ws.Range("K43").Value = 2953.4
$ws.Range("L43").Value = 4650
$ws.Range("M43").Value = -2884.4
$ws.Range("N43").Value = -4788

$ws.Range("H51").Value = 17405.47
$ws.Range("I51").Value = 20983.334
$ws.Range("J51").Value = 15453.909
$ws.Range("K51").Value = 20983.334
$ws.Range("L51").Value = 15453.909
$ws.Range("M51").Value = -20499.334
$ws.Range("N51").Value = -16421.909

$ws.Range("H76").Value = 71380216
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 142756430
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 142756430
$ws.Range("M76").Value = -3685
$ws.Range("N76").Value = -142757060

$ws.Range("H79").Value = 71380216
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 142756430
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 142756430
$ws.Range("M79").Value = -2908
$ws.Range("N79").Value = -142758614

$ws.Range("H107").Value = 22729020
$ws.Range("I107").Value = 13890995
$ws.Range("J107").Value = 62500136
$ws.Range("K107").Value = 13890995
$ws.Range("L107").Value = 62500136
$ws.Range("M107").Value = -13889075
$ws.Range("N107").Value = -62503976

$ws.Range("H111").Value = 3691.4614
$ws.Range("I111").Value = 3118
$ws.Range("J111").Value = 4609
$ws.Range("K111").Value = 9354
$ws.Range("L111").Value = 13827
$ws.Range("M111").Value = -6287
$ws.Range("N111").Value = -19961

$ws.Range("H113").Value = 7489.9614
$ws.Range("I113").Value = 9267.333000000001
$ws.Range("K113").Value = 9267.333000000001
$ws.Range("M113").Value = -6013.333000000001

$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139

$ws.Range("H137").Value = 22794.963
$ws.Range("I137").Value = 32198.107
$ws.Range("K137").Value = 96594.321
$ws.Range("M137").Value = -94044.321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6602.9424
$ws.Range("I32").Value = 4914.6626
$ws.Range("K32").Value = 4914.6626
$ws.Range("M32").Value = -4627.6626

$ws.Range("H61").Value = 3163.361
$ws.Range("I61").Value = 1865.1724
$ws.Range("K61").Value = 1865.1724
$ws.Range("M61").Value = -1653.1724

$ws.Range("H122").Value = 3498.8928
$ws.Range("I122").Value = 2044.3043
$ws.Range("K122").Value = 6132.9129
$ws.Range("M122").Value = -3682.9129

$ws.Range("H135").Value = 57356.75
$ws.Range("J135").Value = 57356.75
$ws.Range("L135").Value = 57356.75
$ws.Range("N135").Value = -67496.75

$ws.Range("H136").Value = 3163.361
$ws.Range("I136").Value = 1865.1724
$ws.Range("K136").Value = 5595.5172
$ws.Range("M136").Value = -3045.5172

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2924.5454
$ws.Range("I105").Value = 2917.1
$ws.Range("K105").Value = 2917.1
$ws.Range("M105").Value = -1170.1

$ws.Range("H134").Value = 2198.6562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1757.8182
$ws.Range("I16").Value = 1637.3334
$ws.Range("J16").Value = 2300
$ws.Range("K16").Value = 1637.3334
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = -1350.3334
$ws.Range("N16").Value = -2874

$ws.Range("H22").Value = 49377.75
$ws.Range("I22").Value = 46005.5
$ws.Range("K22").Value = 46005.5
$ws.Range("M22").Value = -45655.5

$ws.Range("H31").Value = 437176.78
$ws.Range("I31").Value = 1668432.5
$ws.Range("J31").Value = 2615.9412
$ws.Range("K31").Value = 1668432.5
$ws.Range("L31").Value = 2615.9412
$ws.Range("M31").Value = -1668137.5
$ws.Range("N31").Value = -3205.9412

$ws.Range("H34").Value = 437176.78
$ws.Range("I34").Value = 1668432.5
$ws.Range("J34").Value = 2615.9412
$ws.Range("K34").Value = 1668432.5
$ws.Range("L34").Value = 2615.9412
$ws.Range("M34").Value = -1668230.5
$ws.Range("N34").Value = -3019.9412

$ws.Range("H58").Value = 2175.8333
$ws.Range("I58").Value = 1951.24
$ws.Range("J58").Value = 3298.8
$ws.Range("K58").Value = 1951.24
$ws.Range("L58").Value = 3298.8
$ws.Range("M58").Value = -1748.24
$ws.Range("N58").Value = -3704.8

$ws.Range("H99").Value = 848952.2
$ws.Range("I99").Value = 2014050.2
$ws.Range("J99").Value = 16739.285
$ws.Range("K99").Value = 2014050.2
$ws.Range("L99").Value = 16739.285
$ws.Range("M99").Value = -2012552.2
$ws.Range("N99").Value = -19735.285

$ws.Range("H107").Value = 1619.2
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1619.2
$ws.Range("K107").Value = 0
$ws.Range("N107").Value = -5459.2
$ws.Range("M107").ClearContents()

$ws.Range("H113").Value = 1757.8182
$ws.Range("I113").Value = 1637.3334
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1637.3334
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 532.6666
$ws.Range("N113").Value = -6640

$ws.Range("H126").Value = 848952.2
$ws.Range("I126").Value = 2014050.2
$ws.Range("J126").Value = 16739.285
$ws.Range("K126").Value = 6042150.6
$ws.Range("L126").Value = 50217.855
$ws.Range("M126").Value = -6039680.6
$ws.Range("N126").Value = -55157.855

$ws.Range("H134").Value = 2940.2886
$ws.Range("I134").Value = 3350.9524
$ws.Range("J134").Value = 1215.5
$ws.Range("K134").Value = 10052.8572
$ws.Range("L134").Value = 3646.5
$ws.Range("M134").Value = -7517.8572
$ws.Range("N134").Value = -8716.5

$ws.Range("H136").Value = 2175.8333
$ws.Range("I136").Value = 1951.24
$ws.Range("J136").Value = 3298.8
$ws.Range("K136").Value = 5853.72
$ws.Range("L136").Value = 9896.400000000001
$ws.Range("M136").Value = -3303.72
$ws.Range("N136").Value = -14996.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3143.8333
$ws.Range("J137").Value = 2966.5
$ws.Range("L137").Value = 8899.5
$ws.Range("N137").Value = -19099.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.75
$ws.Range("I2").Value = 82.2
$ws.Range("J2").Value = 134.33333
$ws.Range("K2").Value = 82.2
$ws.Range("L2").Value = 134.33333
$ws.Range("M2").Value = 30.8
$ws.Range("N2").Value = -360.33333

$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 10000
$ws.Range("K11").Value = 10000
$ws.Range("M11").Value = -9861

$ws.Range("H97").Value = 2637.2285
$ws.Range("I97").Value = 2134.3635
$ws.Range("K97").Value = 2134.3635
$ws.Range("M97").Value = -1638.3635

$ws.Range("H132").Value = 37565.92
$ws.Range("I132").Value = 37770.406
$ws.Range("J132").Value = 30000
$ws.Range("K132").Value = 113311.218
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -110781.218
$ws.Range("N132").Value = -95060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 31522.592
$ws.Range("I29").Value = 28261.809
$ws.Range("K29").Value = 28261.809
$ws.Range("M29").Value = -27966.809

$ws.Range("H122").Value = 6941.375
$ws.Range("I122").Value = 7286.375
$ws.Range("J122").Value = 6596.375
$ws.Range("K122").Value = 21859.125
$ws.Range("L122").Value = 19789.125
$ws.Range("M122").Value = -19409.125
$ws.Range("N122").Value = -24689.125

$ws.Range("H132").Value = 6120.385
$ws.Range("I132").Value = 4279.4443
$ws.Range("J132").Value = 7698.3335
$ws.Range("K132").Value = 12838.3329
$ws.Range("L132").Value = 23095.0005
$ws.Range("M132").Value = -10308.3329
$ws.Range("N132").Value = -28155.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17868212
$ws.Range("J62").Value = 22738224
$ws.Range("L62").Value = 22738224
$ws.Range("N62").Value = -22739472

$ws.Range("H65").Value = 17868212
$ws.Range("J65").Value = 22738224
$ws.Range("L65").Value = 113691120
$ws.Range("N65").Value = -113697360

$ws.Range("H96").Value = 129499.625
$ws.Range("I96").Value = 338508
$ws.Range("K96").Value = 338508
$ws.Range("M96").Value = -337135

$ws.Range("H107").Value = 1179.6875
$ws.Range("I107").Value = 1063.8572
$ws.Range("J107").Value = 1990.5
$ws.Range("K107").Value = 3191.5716
$ws.Range("L107").Value = 5971.5
$ws.Range("M107").Value = -1271.5716
$ws.Range("N107").Value = -9811.5

$ws.Range("H132").Value = 3603.6924
$ws.Range("I132").Value = 2994.889
$ws.Range("J132").Value = 4973.5
$ws.Range("K132").Value = 8984.667000000001
$ws.Range("L132").Value = 14920.5
$ws.Range("M132").Value = -6454.667000000001
$ws.Range("N132").Value = -19980.5

$ws.Range("H136").Value = 203899.66
$ws.Range("I136").Value = 241671.58
$ws.Range("K136").Value = 725014.74
$ws.Range("M136").Value = -722464.74

$ws.Range("H138").Value = 99998.5
$ws.Range("J138").Value = 99998.5
$ws.Range("L138").Value = 99998.5
$ws.Range("N138").Value = -110278.5
